$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused header cells (C1, D1) and all of rows 5-10,
# which held the PET / Nylon-6_exp2 data that is being dropped.
$ws.Range("C1:D1").Clear()
$ws.Range("A5:D10").Clear()

# Update remaining data rows 2-4 to hold only the PSU series.
$ws.Range("A2").Value = "PSU S1"
$ws.Range("B2").Value = "{0.0: 2494.9877810440958}"
$ws.Range("D2").Clear()

$ws.Range("A3").Value = "PSU S2"
$ws.Range("B3").Value = "{0.0: 2430.1591997439973}"
$ws.Range("D3").Clear()

$ws.Range("A4").Value = "PSU S3"
$ws.Range("B4").Value = "{0.0: 2330.071670255158}"
$ws.Range("D4").Clear()
